$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full org hierarchy data for rows 2-15 (EMPLOYEE.ID, EMPLOYEE name, MANAGER.ID, SALES)
# Order matters: cells are written row-by-row, column-by-column (A,B,C,D) top to bottom
# so that the workbook's shared-strings table is rebuilt in the same order as the
# reference edit.
$data = @(
    @("Ruben.Pruyn2",      "Ruben Pruyn",      "Maira.Roberts3",  70),
    @("Olene.Markiewicz1", "Olene Markiewicz", "Maira.Roberts3",  46),
    @("Kyra.Verra0",       "Kyra Verra",       "Maira.Roberts3",  37),
    @("Maira.Roberts3",    "Maira Roberts",    $null,             80),
    @("Lorna.Hasbell7",    "Lorna Hasbell",    "Kimbra.Agnew4",   90),
    @("Nola.Terstage0",    "Nola Terstage",    "Kimbra.Agnew4",   97),
    @("Earnest.Gwynn5",    "Earnest Gwynn",    "Kimbra.Agnew4",   55),
    @("Kimbra.Agnew4",     "Kimbra Agnew",     $null,             69),
    @("Dede.Waligora9",    "Dede Waligora",    "Ruben.Pruyn2",    91),
    @("Gaynell.Toyota1",   "Gaynell Toyota",   "Ruben.Pruyn2",    97),
    @("Pamela.Harvilla6",  "Pamela Harvilla",  "Ruben.Pruyn2",    57),
    @("Deadra.Ciullo5",    "Deadra Ciullo",    "Lorna.Hasbell7",  69),
    @("Corrinne.Pesch6",   "Corrinne Pesch",   "Lorna.Hasbell7",  49),
    @("Michaele.Trucks6",  "Michaele Trucks",  "Lorna.Hasbell7",  25)
)

$row = 2
foreach ($item in $data) {
    $ws.Cells.Item($row, 1).Value = $item[0]
    $ws.Cells.Item($row, 2).Value = $item[1]
    if ($item[2] -ne $null) {
        $ws.Cells.Item($row, 3).Value = $item[2]
    }
    $ws.Cells.Item($row, 4).Value = $item[3]
    $row++
}
